$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2577
$ws1.Range("F9").Value  = 1260
$ws1.Range("F10").Value = 912
$ws1.Range("F13").Value = 1149
$ws1.Range("F15").Value = 291
$ws1.Range("F17").Value = 736
$ws1.Range("F18").Value = 784
$ws1.Range("F19").Value = 209
$ws1.Range("F20").Value = 500
$ws1.Range("F21").Value = 1125
$ws1.Range("F22").Value = 97
$ws1.Range("F23").Value = 617
$ws1.Range("F28").Value = 685
$ws1.Range("F29").Value = 497
$ws1.Range("F30").Value = 4512
$ws1.Range("F31").Value = 489
$ws1.Range("F35").Value = 158
$ws1.Range("F38").Value = 56
$ws1.Range("F39").Value = 442
$ws1.Range("F45").Value = 136
$ws1.Range("F48").Value = 24

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 190
$ws2.Range("F14").Value = 24
$ws2.Range("F16").Value = 32
$ws2.Range("F22").Value = 2

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 2577
$ws4.Range("F10").Value = 1260
$ws4.Range("F11").Value = 912
$ws4.Range("F14").Value = 1149
$ws4.Range("F15").Value = 291
$ws4.Range("F18").Value = 736
$ws4.Range("F19").Value = 784
$ws4.Range("F20").Value = 209
$ws4.Range("F21").Value = 500
$ws4.Range("F22").Value = 1125
$ws4.Range("F24").Value = 97
$ws4.Range("F25").Value = 617
$ws4.Range("F30").Value = 497
$ws4.Range("F31").Value = 4512
$ws4.Range("F32").Value = 190
$ws4.Range("F33").Value = 489
$ws4.Range("F36").Value = 158
$ws4.Range("F38").Value = 442
$ws4.Range("F39").Value = 24
$ws4.Range("F40").Value = 24
$ws4.Range("F46").Value = 136
